{"js": "// Add select item for Applicable version\n// 1) Update the existing heading + table-cell texts.\n// 2) Append a new \"fg\" Heading2 paragraph followed by a new 4-column\n//    table (header row only) identical in style to the first table's\n//    header row.\n\nconst body = context.document.body;\n\nasync function replaceWholeWord(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait replaceWholeWord(\"c\", \"pcdu\");\nawait replaceWholeWord(\"zdzd\", \"fjvdvdv\");\nawait replaceWholeWord(\"deded\", \"dcdcd\");\nawait replaceWholeWord(\"eedd\", \"cddcdc\");\nawait context.sync();\n\n// Append the new heading paragraph + table via raw OOXML so the new\n// table's header row exactly matches the formatting (bold, red,\n// underlined, sz 24) used by the first table's header row.\nconst newContentOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:pPr><w:pStyle w:val=\"Heading2\"/></w:pPr><w:r><w:t>fg</w:t></w:r></w:p>\n<w:tbl>\n<w:tblPr><w:tblW w:type=\"auto\" w:w=\"0\"/><w:tblLayout w:type=\"fixed\"/><w:tblLook w:firstColumn=\"1\" w:firstRow=\"1\" w:lastColumn=\"0\" w:lastRow=\"0\" w:noHBand=\"0\" w:noVBand=\"1\" w:val=\"04A0\"/></w:tblPr>\n<w:tblGrid><w:gridCol w:w=\"2160\"/><w:gridCol w:w=\"2160\"/><w:gridCol w:w=\"2160\"/><w:gridCol w:w=\"2160\"/></w:tblGrid>\n<w:tr>\n<w:tc><w:tcPr><w:tcW w:type=\"dxa\" w:w=\"2160\"/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:color w:val=\"FF0000\"/><w:sz w:val=\"24\"/><w:u w:val=\"single\"/></w:rPr><w:t>Index de l'etape</w:t></w:r></w:p></w:tc>\n<w:tc><w:tcPr><w:tcW w:type=\"dxa\" w:w=\"2160\"/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:color w:val=\"FF0000\"/><w:sz w:val=\"24\"/><w:u w:val=\"single\"/></w:rPr><w:t>Nom de l'etape</w:t></w:r></w:p></w:tc>\n<w:tc><w:tcPr><w:tcW w:type=\"dxa\" w:w=\"2160\"/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:color w:val=\"FF0000\"/><w:sz w:val=\"24\"/><w:u w:val=\"single\"/></w:rPr><w:t>Description de l'\u00e9tape</w:t></w:r></w:p></w:tc>\n<w:tc><w:tcPr><w:tcW w:type=\"dxa\" w:w=\"2160\"/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:color w:val=\"FF0000\"/><w:sz w:val=\"24\"/><w:u w:val=\"single\"/></w:rPr><w:t>R\u00e9sultat Attendu</w:t></w:r></w:p></w:tc>\n</w:tr>\n</w:tbl>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nbody.insertOoxml(newContentOoxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Add select item for Applicable version\n$d = $word.ActiveDocument\n\n# --- 1) Update the existing heading + table-cell texts. ---\nfunction Replace-WholeWord($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-WholeWord \"c\" \"pcdu\"\nReplace-WholeWord \"zdzd\" \"fjvdvdv\"\nReplace-WholeWord \"deded\" \"dcdcd\"\nReplace-WholeWord \"eedd\" \"cddcdc\"\n\n# --- 2) Append a new \"fg\" Heading2 paragraph followed by a new 4-column ---\n#        table (header row only), matching the first table's header-row\n#        formatting (bold, red, underlined, sz 24).\n#\n# A bare Paragraphs.Add() first gives us a true body-level (not\n# in-table-cell) insertion point at the very end of the document; we then\n# replace that placeholder paragraph's range with exact OOXML so the\n# emitted markup matches the original table's formatting precisely.\n$placeholder = $d.Paragraphs.Add()\n$placeholderRange = $placeholder.Range\n\n$ooxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:pPr><w:pStyle w:val=\"Heading2\"/></w:pPr><w:r><w:t>fg</w:t></w:r></w:p>\n<w:tbl>\n<w:tblPr><w:tblW w:type=\"auto\" w:w=\"0\"/><w:tblLayout w:type=\"fixed\"/><w:tblLook w:firstColumn=\"1\" w:firstRow=\"1\" w:lastColumn=\"0\" w:lastRow=\"0\" w:noHBand=\"0\" w:noVBand=\"1\" w:val=\"04A0\"/></w:tblPr>\n<w:tblGrid><w:gridCol w:w=\"2160\"/><w:gridCol w:w=\"2160\"/><w:gridCol w:w=\"2160\"/><w:gridCol w:w=\"2160\"/></w:tblGrid>\n<w:tr>\n<w:tc><w:tcPr><w:tcW w:type=\"dxa\" w:w=\"2160\"/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:color w:val=\"FF0000\"/><w:sz w:val=\"24\"/><w:u w:val=\"single\"/></w:rPr><w:t>Index de l'etape</w:t></w:r></w:p></w:tc>\n<w:tc><w:tcPr><w:tcW w:type=\"dxa\" w:w=\"2160\"/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:color w:val=\"FF0000\"/><w:sz w:val=\"24\"/><w:u w:val=\"single\"/></w:rPr><w:t>Nom de l'etape</w:t></w:r></w:p></w:tc>\n<w:tc><w:tcPr><w:tcW w:type=\"dxa\" w:w=\"2160\"/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:color w:val=\"FF0000\"/><w:sz w:val=\"24\"/><w:u w:val=\"single\"/></w:rPr><w:t>Description de l'\u00e9tape</w:t></w:r></w:p></w:tc>\n<w:tc><w:tcPr><w:tcW w:type=\"dxa\" w:w=\"2160\"/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:color w:val=\"FF0000\"/><w:sz w:val=\"24\"/><w:u w:val=\"single\"/></w:rPr><w:t>R\u00e9sultat Attendu</w:t></w:r></w:p></w:tc>\n</w:tr>\n</w:tbl>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n'@\n\n$placeholderRange.InsertXML($ooxml)\n"}
